$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new data row above row 25 (pushes "VITAMOUNT..." and everything
#    below it, including the totals & footer rows, down by one row).
# ---------------------------------------------------------------------------
$ws.Rows("25:25").Insert()

# The freshly-inserted row loses its row height and (because Excel only
# carries borders that are shared by the rows on *both* sides) its bottom
# border. Restore both so the row matches its neighbours (rows 24/26..32
# all use the same "data row" look).
$ws.Rows("25:25").RowHeight = 24.75

$newRowBorder = $ws.Range("A25:Q25").Borders.Item(9)   # xlEdgeBottom
$refBorder = $ws.Range("A26").Borders.Item(9)
$newRowBorder.LineStyle = $refBorder.LineStyle
$newRowBorder.Weight = $refBorder.Weight
$newRowBorder.Color = $refBorder.Color

# ---------------------------------------------------------------------------
# 2. Populate the new row with the new item: ULCEZOLE 60MG 20 DEL. REL.
#    CAPSULES.
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = 19
$ws.Range("C25").Value = "ULCEZOLE 60MG 20 DEL. REL. CAPSULES"
$ws.Range("H25").Value = "0:1"
$ws.Range("L25").Value = "1"
$ws.Range("N25").Value = "198.00"
$ws.Range("P25").Value = "198.0000"
$ws.Range("Q25").Value = "1:0"

# ---------------------------------------------------------------------------
# 3. Renumber the "#" column for every row that shifted down one position
#    (their literal counters need to stay sequential: 20, 21, 22, ...).
# ---------------------------------------------------------------------------
for ($r = 26; $r -le 32; $r++) {
    $ws.Range("A$r").Value = $r - 6
}

# ---------------------------------------------------------------------------
# 4. Update the grand-total cell (now on row 33) to include the new item's
#    price.
# ---------------------------------------------------------------------------
$ws.Range("P33").Value = 1597.01

# ---------------------------------------------------------------------------
# 5. Refresh the generated-on timestamp in the footer (now on row 34).
# ---------------------------------------------------------------------------
$ws.Range("A34").Value = "Wednesday, 28 May, 2025 3:20 PM"
